$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (分值 - score per question) changes from 5 to 8.33 for every
# answered row, and several answer letters in column B are revised.
# Row 10 (question 9) is cleared out entirely (no answer recorded yet),
# while row 14 (question 13) is newly populated with a result.

$ws.Range("B2").Value = "C"
$ws.Range("C2").Value = 8.33

$ws.Range("B3").Value = "C"
$ws.Range("C3").Value = 8.33

$ws.Range("B4").Value = "D"
$ws.Range("C4").Value = 8.33

$ws.Range("B5").Value = "B"
$ws.Range("C5").Value = 8.33

$ws.Range("B6").Value = "A"
$ws.Range("C6").Value = 8.33

$ws.Range("B7").Value = "D"
$ws.Range("C7").Value = 8.33

$ws.Range("B8").Value = "B"
$ws.Range("C8").Value = 8.33

$ws.Range("C9").Value = 8.33

$ws.Range("B10").ClearContents()
$ws.Range("C10").ClearContents()
$ws.Range("D10").ClearContents()

$ws.Range("B11").Value = "C"
$ws.Range("C11").Value = 8.33

$ws.Range("C12").Value = 8.33

$ws.Range("B13").Value = "C"
$ws.Range("C13").Value = 8.33

$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "B"
$ws.Range("C14").Value = 8.33
$ws.Range("D14").Value = 0

$ws.Range("D16").Select() | Out-Null
